$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 used to be "Nymphali" (and the old rows 4-5, Nymphali/Nocatli, are
# dropped); what used to be row 6 ("Giratina Holo") becomes the new row 4.
$ws.Range("A4").Value = "Giratina Holo"
$ws.Range("B4").Value = "10/127"
$ws.Range("C4").Value = "Platine"
$ws.Range("D4").Value = "FR"
$ws.Range("E4").Value = "https://www.cardmarket.com/fr/Pokemon/Products/Singles/Platinum/Giratina-Lv63-PL10"

# Seed rows 5-11 with row 4's formatting (general alignment text columns,
# right-aligned Serie column) before filling in the new values, so Excel
# doesn't fall back to the plain column default style for the newly used
# rows.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E11").PasteSpecial(-4122)

# Row 5: Giratina VSTAR
$ws.Range("A5").Value = "Giratina VSTAR"
$ws.Range("B5").Value = "261/172"
$ws.Range("C5").Value = "S12A"
$ws.Range("D5").Value = "JP"
$ws.Range("E5").Value = "https://www.cardmarket.com/fr/Pokemon/Products/Singles/VSTAR-Universe/Giratina-VSTAR-V2-s12a261"

# Row 6: Energie Tenebre
$ws.Range("A6").Value = "Energie Tênebre"
$ws.Range("B6").Value = "257/172"
$ws.Range("C6").Value = "S12A"
$ws.Range("D6").Value = "JP"
$ws.Range("E6").Value = "https://www.cardmarket.com/fr/Pokemon/Products/Singles/VSTAR-Universe/Darkness-Energy-s12a257"

# Row 7: Energie Acier
$ws.Range("A7").Value = "Energie Acier"
$ws.Range("B7").Value = "258/172"
$ws.Range("C7").Value = "S12A"
$ws.Range("D7").Value = "JP"
$ws.Range("E7").Value = "https://www.cardmarket.com/fr/Pokemon/Products/Singles/VSTAR-Universe/Metal-Energy-s12a258"

# Row 8: Energie Feu
$ws.Range("A8").Value = "Energie Feu"
$ws.Range("B8").Value = "252/172"
$ws.Range("C8").Value = "S12A"
$ws.Range("D8").Value = "JP"
$ws.Range("E8").Value = "https://www.cardmarket.com/fr/Pokemon/Products/Singles/VSTAR-Universe/Fire-Energy-s12a252"

# Row 9: Regigigas VSTAR
$ws.Range("A9").Value = "Regigigas VSTAR"
$ws.Range("B9").Value = "233/172"
$ws.Range("C9").Value = "S12A"
$ws.Range("D9").Value = "JP"
$ws.Range("E9").Value = "https://www.cardmarket.com/fr/Pokemon/Products/Singles/VSTAR-Universe/Regigigas-VSTAR-V2-s12a233"

# Row 10: Entei V
$ws.Range("A10").Value = "Entei V"
$ws.Range("B10").Value = "213/172"
$ws.Range("C10").Value = "S12A"
$ws.Range("D10").Value = "JP"
$ws.Range("E10").Value = "https://www.cardmarket.com/fr/Pokemon/Products/Singles/VSTAR-Universe/Entei-V-V2-s12a213"

# Row 11: Charisme de Giovanni (C11 is a plain number, not text)
$ws.Range("A11").Value = "Charisme de Giovanni"
$ws.Range("B11").Value = "204/165"
$ws.Range("C11").Value = 151
$ws.Range("D11").Value = "FR"
$ws.Range("E11").Value = "https://www.cardmarket.com/fr/Pokemon/Products/Singles/151/Giovannis-Charisma-V3-MEW204"

# Apply a numeric (thousands separator) format to the "Serie" column so the
# new plain-number entry (C11 = 151) renders like the rest of the column.
$ws.Range("C1:C11").NumberFormat = "#,##0"

# Column F is no longer used; delete it entirely.
$ws.Columns("F").Delete()
